# Update Name of Algo
# Applies updated numeric values produced by a re-run of the RandomForest
# imputation algorithm to the corresponding cells on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A3"   = -22.10169999999999
    "C6"   = -12.38779999999999
    "D10"  = -8.530400000000004
    "A14"  = -22.06799999999999
    "D15"  = -8.5829
    "C18"  = -11.11989999999999
    "C19"  = -12.61020000000001
    "A21"  = -19.82609999999999
    "D21"  = -7.713299999999997
    "B22"  = 8.631300000000005
    "D22"  = -8.328600000000005
    "A23"  = -20.09609999999997
    "B24"  = 6.001
    "D24"  = -7.062999999999997
    "A25"  = -21.94399999999999
    "A26"  = -21.03509999999996
    "B28"  = 6.000900000000001
    "A29"  = -20.86209999999998
    "D33"  = -8.269399999999997
    "B36"  = 9.301600000000006
    "C44"  = -12.80290000000001
    "B45"  = 4.663400000000005
    "D46"  = -7.981799999999997
    "C47"  = -12.30550000000001
    "B48"  = 5.221500000000005
    "B49"  = 5.5336
    "D49"  = -8.271199999999999
    "C51"  = -10.8494
    "B52"  = 5.310899999999997
    "A53"  = -21.76280000000001
    "B53"  = 5.955799999999996
    "B54"  = 4.996900000000002
    "C55"  = -13.401
    "D56"  = -8.531799999999997
    "A57"  = -22.02539999999999
    "C57"  = -12.7922
    "A59"  = -22.5653
    "D61"  = -8.339099999999997
    "C64"  = -10.43899999999999
    "D66"  = -7.006299999999997
    "A69"  = -21.656
    "B70"  = 5.026900000000001
    "D74"  = -8.336700000000009
    "D77"  = -6.437799999999997
    "A79"  = -20.27330000000001
    "C80"  = -13.1625
    "A83"  = -21.69200000000001
    "B86"  = 5.068300000000002
    "B87"  = 5.437099999999996
    "D87"  = -8.508599999999999
    "D88"  = -7.624199999999997
    "B89"  = 4.548299999999998
    "A91"  = -20.24869999999997
    "C92"  = -10.28770000000001
    "A93"  = -21.29400000000002
    "C94"  = -10.6689
    "C96"  = -10.27260000000001
    "D100" = -8.122999999999998
    "B101" = 6.1287
    "C101" = -12.42719999999999
    "A103" = -21.8485
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
